$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new X-1P-User header alongside the existing Content-Type header in
# the HEADERS column (F) for every test row that sends that header.
$newHeader = "X-1P-User=(SYS_USER1)||Content-Type=application/json"
$ws.Range("F3").Value = $newHeader
$ws.Range("F5").Value = $newHeader
$ws.Range("F7").Value = $newHeader
$ws.Range("F8").Value = $newHeader
$ws.Range("F9").Value = $newHeader
$ws.Range("F10").Value = $newHeader

# The HEADERS column needs to be widened now that its values are longer.
$ws.Columns("F").ColumnWidth = 51.666666666666664

# Clear out the now-unused STATUS column (L) results.
$ws.Range("L2:L10").ClearContents()

# Leave the selection on the cleared range, matching the state the workbook
# was left in after the edit.
$ws.Range("L2:L10").Select()
